$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header "VIMMP_DEF" in column F, matching the style of the other headers
$ws.Range("F1").Value = "VIMMP_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Fill new column F with "[]" for the 9 data rows (no special style, like columns C/E)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
